$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.0563737760172628
$ws.Range("C2").Value = 0.0563737760172628
$ws.Range("D2").Value = 0.0409130078879325
$ws.Range("E2").Value = 0.000282451268113067
$ws.Range("F2").Value = 0.8338

# Row 3
$ws.Range("B3").Value = 0.761313018832472
$ws.Range("C3").Value = 0.761313018832472
$ws.Range("D3").Value = 0.552519411421733
$ws.Range("E3").Value = 0.00381443009129584
$ws.Range("F3").Value = 0.4595

# Row 4
$ws.Range("B4").Value = 0.353238455572906
$ws.Range("C4").Value = 0.353238455572906
$ws.Range("D4").Value = 0.256361179615676
$ws.Range("F4").Value = 0.6183
